$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting username/password to B/C
$ws.Range("A1").EntireColumn.Insert()

# New "testType" column values
$ws.Range("A1").Value = "testType"
$ws.Range("A2").Value = "happyPath"
$ws.Range("A3").Value = "errorPath"
$ws.Range("A4").Value = "errorPath"
$ws.Range("A5").Value = "errorPath2"

# New row 5: reuse FirstTestLogin_12 with an incorrect password
$ws.Range("B5").Value = "FirstTestLogin_12"
$ws.Range("C5").Value = "incorrect"

# Selection / active cell as recorded in the saved file
$ws.Range("B2").Select()
